# "testing with 2004 as representative"
# Update the "Coupling Parameters" sheet: flip the two demand/profile
# representative-year switches off, change the future year to 2004,
# remove the now-obsolete "not implemented" helper cell in B46, shrink
# the conditional formatting range accordingly, and move the active
# selection to B31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Demand same as representative year? -> No
$ws.Range("B28").Value = $false

# Actual year and investment year are representative year? -> No
$ws.Range("B29").Value = $false

# Representative / future year to test with
$ws.Range("B30").Value = 2004

# Remove the obsolete helper formula in B46 (row keeps its label in A46)
$ws.Range("B46").ClearContents()

# The conditional formatting that used to start at B46 now starts at B47
$fcs = $ws.Range("B46:B51").FormatConditions
$fc = $fcs.Item(1)
$fc.ModifyAppliesToRange($ws.Range("B47:B51"))

# Leave the selection on B31, matching the saved view state
$ws.Range("B31").Select()
